$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet stores every value (even "4", "0", dates, etc.) as literal TEXT
# (shared strings), not as numbers/dates. Typing values straight into a Range
# via .Value would let Excel "smart type" them into numbers/dates, which would
# not match the source data. Every value we need already exists somewhere on
# the sheet, so instead we relocate/duplicate existing cells with Copy +
# PasteSpecial(xlPasteValues = -4163), which preserves the literal text type,
# the shared-string reuse, and each destination cell's own style/format.
# ---------------------------------------------------------------------------

$xlPasteValues = -4163

# --- Phase 1: snapshot every distinct source cell we will need into a scratch column (far off-sheet) ---
$ws.Range("B1").Copy() | Out-Null
$ws.Range("Z2000").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("C1").Copy() | Out-Null
$ws.Range("Z2001").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("Z2002").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A3").Copy() | Out-Null
$ws.Range("Z2003").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("Z2004").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A4").Copy() | Out-Null
$ws.Range("Z2005").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("Z2006").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A5").Copy() | Out-Null
$ws.Range("Z2007").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("B5").Copy() | Out-Null
$ws.Range("Z2008").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A6").Copy() | Out-Null
$ws.Range("Z2009").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("B6").Copy() | Out-Null
$ws.Range("Z2010").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A7").Copy() | Out-Null
$ws.Range("Z2011").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("B7").Copy() | Out-Null
$ws.Range("Z2012").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A8").Copy() | Out-Null
$ws.Range("Z2013").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("B8").Copy() | Out-Null
$ws.Range("Z2014").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("Z2015").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("B9").Copy() | Out-Null
$ws.Range("Z2016").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A10").Copy() | Out-Null
$ws.Range("Z2017").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("B13").Copy() | Out-Null
$ws.Range("Z2018").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A11").Copy() | Out-Null
$ws.Range("Z2019").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A12").Copy() | Out-Null
$ws.Range("Z2020").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("Z2021").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A16").Copy() | Out-Null
$ws.Range("Z2022").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A17").Copy() | Out-Null
$ws.Range("Z2023").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A18").Copy() | Out-Null
$ws.Range("Z2024").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A19").Copy() | Out-Null
$ws.Range("Z2025").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A20").Copy() | Out-Null
$ws.Range("Z2026").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("B14").Copy() | Out-Null
$ws.Range("Z2027").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A21").Copy() | Out-Null
$ws.Range("Z2028").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("B20").Copy() | Out-Null
$ws.Range("Z2029").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A22").Copy() | Out-Null
$ws.Range("Z2030").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("B21").Copy() | Out-Null
$ws.Range("Z2031").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A23").Copy() | Out-Null
$ws.Range("Z2032").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("B22").Copy() | Out-Null
$ws.Range("Z2033").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("A24").Copy() | Out-Null
$ws.Range("Z2034").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("B25").Copy() | Out-Null
$ws.Range("Z2035").PasteSpecial($xlPasteValues) | Out-Null

# --- Phase 2: write the snapshots into their final destination cells ---
$ws.Range("Z2000").Copy() | Out-Null
$ws.Range("B1").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2001").Copy() | Out-Null
$ws.Range("C1").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2002").Copy() | Out-Null
$ws.Range("B2").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2002").Copy() | Out-Null
$ws.Range("C2").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2003").Copy() | Out-Null
$ws.Range("A3").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2004").Copy() | Out-Null
$ws.Range("B3").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2004").Copy() | Out-Null
$ws.Range("C3").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2005").Copy() | Out-Null
$ws.Range("A4").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2006").Copy() | Out-Null
$ws.Range("B4").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2006").Copy() | Out-Null
$ws.Range("C4").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2007").Copy() | Out-Null
$ws.Range("A5").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2008").Copy() | Out-Null
$ws.Range("B5").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2008").Copy() | Out-Null
$ws.Range("C5").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2009").Copy() | Out-Null
$ws.Range("A6").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2010").Copy() | Out-Null
$ws.Range("B6").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2010").Copy() | Out-Null
$ws.Range("C6").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2011").Copy() | Out-Null
$ws.Range("A7").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2012").Copy() | Out-Null
$ws.Range("B7").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2012").Copy() | Out-Null
$ws.Range("C7").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2013").Copy() | Out-Null
$ws.Range("A8").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2014").Copy() | Out-Null
$ws.Range("B8").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2014").Copy() | Out-Null
$ws.Range("C8").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2015").Copy() | Out-Null
$ws.Range("A9").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2016").Copy() | Out-Null
$ws.Range("B9").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2016").Copy() | Out-Null
$ws.Range("C9").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2017").Copy() | Out-Null
$ws.Range("A10").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2018").Copy() | Out-Null
$ws.Range("B10").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2018").Copy() | Out-Null
$ws.Range("C10").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2019").Copy() | Out-Null
$ws.Range("A11").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2020").Copy() | Out-Null
$ws.Range("A12").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2021").Copy() | Out-Null
$ws.Range("A13").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2014").Copy() | Out-Null
$ws.Range("B13").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2014").Copy() | Out-Null
$ws.Range("C13").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2022").Copy() | Out-Null
$ws.Range("A14").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2023").Copy() | Out-Null
$ws.Range("A15").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2018").Copy() | Out-Null
$ws.Range("B15").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2018").Copy() | Out-Null
$ws.Range("C15").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2024").Copy() | Out-Null
$ws.Range("A16").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2025").Copy() | Out-Null
$ws.Range("A17").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2026").Copy() | Out-Null
$ws.Range("A18").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2027").Copy() | Out-Null
$ws.Range("B18").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2027").Copy() | Out-Null
$ws.Range("C18").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2028").Copy() | Out-Null
$ws.Range("A19").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2029").Copy() | Out-Null
$ws.Range("B19").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2029").Copy() | Out-Null
$ws.Range("C19").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2030").Copy() | Out-Null
$ws.Range("A20").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2031").Copy() | Out-Null
$ws.Range("B20").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2031").Copy() | Out-Null
$ws.Range("C20").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2032").Copy() | Out-Null
$ws.Range("A21").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2033").Copy() | Out-Null
$ws.Range("B21").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2033").Copy() | Out-Null
$ws.Range("C21").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2034").Copy() | Out-Null
$ws.Range("A22").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2035").Copy() | Out-Null
$ws.Range("B23").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("Z2035").Copy() | Out-Null
$ws.Range("C23").PasteSpecial($xlPasteValues) | Out-Null

$ws.Application.CutCopyMode = $false

# --- Phase 3: remove the scratch column ---
$ws.Range("Z2000:Z2035").ClearContents() | Out-Null

# --- Phase 4: the sheet shrank from 25 data rows to 23; drop the old trailing rows 24-25 ---
$ws.Range("A24:A25").EntireRow.Delete() | Out-Null

# --- Phase 5: fix up row heights to match the new layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
